# refactor: Remove orçamento, finanças e indicadores de saúde do projeto
#
# The "Orçamento" column (F) in the PROCESSO header/data block (rows 2-3)
# is removed, and the "Descrição" column (G) is shifted left into its
# place, so F now holds the description text and G is emptied out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (header): drop "Orçamento" in F2, shift "Descrição" from G2 into F2
$ws.Range("F2").Value2 = $ws.Range("G2").Value2

# Row 3 (data): drop "R$ 15.000" in F3, shift description text from G3 into F3
$ws.Range("F3").Value2 = $ws.Range("G3").Value2

# Remove the now-duplicated G2/G3 cells entirely (Clear, not just
# ClearContents, so the cell elements disappear rather than remaining as
# empty styled placeholders). Other rows' G column must stay untouched.
$ws.Range("G2").Clear()
$ws.Range("G3").Clear()
